# Add the new game title as a new row at the bottom of the list (used to
# set up a test case for the table filter feature, sitting between the
# "adding" and "installing" games flows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "Meant To Fail"

# Leave the selection where it was left after typing the new entry.
$ws.Cells.Item($newRow, 3).Select() | Out-Null
